# Update NATMI Wnt6-Fzd7 sheet with newly computed TPM-based values.
# - Rows 2-9 (the FAPs- and Resolving-Mac-sending groups) get refreshed
#   numeric results (columns G-T) from the new TPM run; column D's
#   "Resolving-Mac"/"MuSCs" split also shifts for rows 4/5 and 8/9.
# - The old "MuSCs"-sending block (rows 6-9) is replaced by the
#   "Resolving-Mac"-sending block, and the trailing "Resolving-Mac"-sending
#   rows (old rows 10-13) are dropped entirely, shrinking the sheet from
#   A1:T13 to A1:T9.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the last four data rows (old rows 10-13) - the sheet shrinks to A1:T9.
$ws.Range("A10:T13").Delete()

function Set-DataRow {
    param($row, $a, $b, $c, $d, $e, $f, $g, $h, $i, $j, $k, $l, $m, $n, $o, $p, $q, $r, $s, $t)

    $ws.Cells.Item($row, 1).Value = $a
    $ws.Cells.Item($row, 2).Value = $b
    $ws.Cells.Item($row, 3).Value = $c
    $ws.Cells.Item($row, 4).Value = $d
    $ws.Cells.Item($row, 5).Value = $e
    $ws.Cells.Item($row, 6).Value = $f
    $ws.Cells.Item($row, 7).Value = $g
    $ws.Cells.Item($row, 8).Value = $h
    $ws.Cells.Item($row, 9).Value = $i
    $ws.Cells.Item($row, 10).Value = $j
    $ws.Cells.Item($row, 11).Value = $k
    $ws.Cells.Item($row, 12).Value = $l
    $ws.Cells.Item($row, 13).Value = $m
    $ws.Cells.Item($row, 14).Value = $n
    $ws.Cells.Item($row, 15).Value = $o
    $ws.Cells.Item($row, 16).Value = $p
    $ws.Cells.Item($row, 17).Value = $q
    $ws.Cells.Item($row, 18).Value = $r
    $ws.Cells.Item($row, 19).Value = $s
    $ws.Cells.Item($row, 20).Value = $t
}

Set-DataRow 2 "FAPs" "Wnt6" "Fzd7" "ECs" `
    1 0.3333333333333333 0.1008403333333333 0.302521 `
    0.5164422077268048 0.5164422077268047 3 1 `
    1.123319 3.369957 0.05053686506648315 0.05053686506648315 `
    0.1132758623996667 1.019482761597 0.0260993701665262 0.02609937016652619

Set-DataRow 3 "FAPs" "Wnt6" "Fzd7" "FAPs" `
    1 0.3333333333333333 0.1008403333333333 0.302521 `
    0.5164422077268048 0.5164422077268047 3 1 `
    11.20764866666667 33.622946 0.5042195746532222 0.5042195746532223 `
    1.130183027429555 10.171647246866 0.2604002703129806 0.2604002703129806

Set-DataRow 4 "FAPs" "Wnt6" "Fzd7" "MuSCs" `
    1 0.3333333333333333 0.1008403333333333 0.302521 `
    0.5164422077268048 0.5164422077268047 3 1 `
    4.958620666666667 14.875862 0.2230827962023326 0.2230827962023326 `
    0.5000289609002223 4.500260648102 0.1152093717766015 0.1152093717766015

Set-DataRow 5 "FAPs" "Wnt6" "Fzd7" "Resolving-Mac" `
    1 0.3333333333333333 0.1008403333333333 0.302521 `
    0.5164422077268048 0.5164422077268047 3 1 `
    4.938126 14.814378 0.222160764077962 0.222160764077962 `
    0.497962271882 4.481660446938 0.1147331954706965 0.1147331954706965

Set-DataRow 6 "Resolving-Mac" "Wnt6" "Fzd7" "ECs" `
    1 0.3333333333333333 0.09441933333333334 0.283258 `
    0.4835577922731952 0.4835577922731952 3 1 `
    1.123319 3.369957 0.05053686506648315 0.05053686506648315 `
    0.1060630311006667 0.9545672799060001 0.02443749489995695 0.02443749489995695

Set-DataRow 7 "Resolving-Mac" "Wnt6" "Fzd7" "FAPs" `
    1 0.3333333333333333 0.09441933333333334 0.283258 `
    0.4835577922731952 0.4835577922731952 3 1 `
    11.20764866666667 33.622946 0.5042195746532222 0.5042195746532223 `
    1.058218715340889 9.523968438068 0.2438193043402417 0.2438193043402417

Set-DataRow 8 "Resolving-Mac" "Wnt6" "Fzd7" "MuSCs" `
    1 0.3333333333333333 0.09441933333333334 0.283258 `
    0.4835577922731952 0.4835577922731952 3 1 `
    4.958620666666667 14.875862 0.2230827962023326 0.2230827962023326 `
    0.4681896575995557 4.213706918396 0.1078734244257311 0.1078734244257311

Set-DataRow 9 "Resolving-Mac" "Wnt6" "Fzd7" "Resolving-Mac" `
    1 0.3333333333333333 0.09441933333333334 0.283258 `
    0.4835577922731952 0.4835577922731952 3 1 `
    4.938126 14.814378 0.222160764077962 0.222160764077962 `
    0.4662545648360001 4.196291083524001 0.1074275686072655 0.1074275686072655

Write-Output "Wnt6-Fzd7 sheet refreshed with new TPM values; rows trimmed to A1:T9"
